$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Punti-storia:" placeholder cells with actual story-point numbers
$ws.Range("K3").Value = "Punti-storia: 2"
$ws.Range("E3").Value = "Punti-storia: 1"

$ws.Range("E15").Value = "Punti-storia: 2"

$ws.Range("E27").Value = "Punti-storia: 2"
$ws.Range("K27").Value = "Punti-storia: 1"

$ws.Range("E39").Value = "Punti-storia: 1"
$ws.Range("K39").Value = "Punti-storia: 2"

$ws.Range("C51").Value = "Priorità: 2"
$ws.Range("E51").Value = "Punti-storia: 2"
$ws.Range("I51").Value = "Priorità: 1"
$ws.Range("K51").Value = "Punti-storia: 3"

# Sheet view: drop the saved scroll position, move selection to O5
$ws.Activate()
$ws.Range("O5").Select()

# Window layout (maximized-ish geometry seen in the target file)
$excel.Left = -108
$excel.Top = -108
$excel.Width = 23256
$excel.Height = 12456

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
